$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Selectors", "Box Model", "Flexbox" and "Grid" (CSS topic rows 12-15) as completed.
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 1

# Move the visible selection / scroll position to E16, matching where the
# user continued working next.
$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollRow = 6
